$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92:125 down to 93:126
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new record's data
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 44726
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112052
$ws.Cells.Item(92, 7).Value = "Albahaca"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 90
$ws.Cells.Item(92, 11).Value = 5000
$ws.Cells.Item(92, 12).Value = 5000
$ws.Cells.Item(92, 13).Value = 5000
$ws.Cells.Item(92, 14).Value = "`$/docena de matas"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 833
$ws.Cells.Item(92, 17).Value = 6
$ws.Cells.Item(92, 18).Value = "Hortaliza"
